# Fruta / hortaliza, semanal
# Weekly refresh of the price-history rows: dates (col D), volume (col M),
# min/max/weighted-avg prices (cols N/O/P) and price-per-kg (col S) are
# updated to reflect the latest weekly data pull. Only the rows that
# actually changed values are touched; rows 4, 7 and 10 are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44357
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 725

# Row 3
$ws.Range("D3").Value = 44893
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21625
$ws.Range("S3").Value = 1081

# Row 5
$ws.Range("D5").Value = 44792
$ws.Range("M5").Value = 100
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 1075

# Row 6
$ws.Range("D6").Value = 44761
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20500
$ws.Range("S6").Value = 1025

# Row 8
$ws.Range("D8").Value = 44708
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("S8").Value = 1025

# Row 9
$ws.Range("D9").Value = 44320
$ws.Range("M9").Value = 80

# Row 11
$ws.Range("D11").Value = 44533
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 17000
$ws.Range("P11").Value = 16500
$ws.Range("S11").Value = 825
